$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write the numeric columns (B:T) for the new rows 511-517 ---
$ws.Range("B511").Value = 510
$ws.Range("C511").Value = 20396
$ws.Range("D511").Value = 581
$ws.Range("E511").Value = 89
$ws.Range("F511").Value = 0.02848597764
$ws.Range("G511").Value = 19667
$ws.Range("H511").Value = 27040
$ws.Range("I511").Value = 47436
$ws.Range("J511").Value = 21
$ws.Range("K511").Value = 1
$ws.Range("L511").Value = 21
$ws.Range("M511").Value = 43
$ws.Range("N511").Value = 0
$ws.Range("O511").Value = 43
$ws.Range("P511").Value = 8
$ws.Range("Q511").Value = 18
$ws.Range("R511").Value = 18
$ws.Range("S511").Value = 61
$ws.Range("T511").Value = 73

$ws.Range("B512").Value = 511
$ws.Range("C512").Value = 20396
$ws.Range("D512").Value = 581
$ws.Range("E512").Value = 89
$ws.Range("F512").Value = 0.02848597764
$ws.Range("G512").Value = 19667
$ws.Range("H512").Value = 27040
$ws.Range("I512").Value = 47436
$ws.Range("J512").Value = 0
$ws.Range("K512").Value = 0
$ws.Range("L512").Value = 0
$ws.Range("M512").Value = 43
$ws.Range("N512").Value = 0
$ws.Range("O512").Value = 43
$ws.Range("P512").Value = 8
$ws.Range("Q512").Value = 18
$ws.Range("R512").Value = 18
$ws.Range("S512").Value = 61
$ws.Range("T512").Value = 74

$ws.Range("B513").Value = 512
$ws.Range("C513").Value = 20439
$ws.Range("D513").Value = 581
$ws.Range("E513").Value = 81
$ws.Range("F513").Value = 0.02842604824
$ws.Range("G513").Value = 19718
$ws.Range("H513").Value = 27040
$ws.Range("I513").Value = 47479
$ws.Range("J513").Value = 43
$ws.Range("K513").Value = 0
$ws.Range("L513").Value = 43
$ws.Range("M513").Value = 80
$ws.Range("N513").Value = 0
$ws.Range("O513").Value = 80
$ws.Range("P513").Value = 7
$ws.Range("Q513").Value = 18
$ws.Range("R513").Value = 17
$ws.Range("S513").Value = 61
$ws.Range("T513").Value = 74

$ws.Range("B514").Value = 513
$ws.Range("C514").Value = 20459
$ws.Range("D514").Value = 581
$ws.Range("E514").Value = 91
$ws.Range("F514").Value = 0.02839825993
$ws.Range("G514").Value = 19728
$ws.Range("H514").Value = 27040
$ws.Range("I514").Value = 47499
$ws.Range("J514").Value = 20
$ws.Range("K514").Value = 0
$ws.Range("L514").Value = 20
$ws.Range("M514").Value = 88
$ws.Range("N514").Value = 0
$ws.Range("O514").Value = 88
$ws.Range("P514").Value = 8
$ws.Range("Q514").Value = 18
$ws.Range("R514").Value = 19
$ws.Range("S514").Value = 61
$ws.Range("T514").Value = 74

$ws.Range("B515").Value = 514
$ws.Range("C515").Value = 20470
$ws.Range("D515").Value = 581
$ws.Range("E515").Value = 94
$ws.Range("F515").Value = 0.02838299951
$ws.Range("G515").Value = 19736
$ws.Range("H515").Value = 27040
$ws.Range("I515").Value = 47510
$ws.Range("J515").Value = 11
$ws.Range("K515").Value = 0
$ws.Range("L515").Value = 11
$ws.Range("M515").Value = 66
$ws.Range("N515").Value = 0
$ws.Range("O515").Value = 66
$ws.Range("P515").Value = 6
$ws.Range("Q515").Value = 17
$ws.Range("R515").Value = 19
$ws.Range("S515").Value = 61
$ws.Range("T515").Value = 74

$ws.Range("B516").Value = 515
$ws.Range("C516").Value = 20471
$ws.Range("D516").Value = 581
$ws.Range("E516").Value = 86
$ws.Range("F516").Value = 0.02838161301
$ws.Range("G516").Value = 19745
$ws.Range("H516").Value = 27040
$ws.Range("I516").Value = 47511
$ws.Range("J516").Value = 1
$ws.Range("K516").Value = 0
$ws.Range("L516").Value = 1
$ws.Range("M516").Value = 66
$ws.Range("N516").Value = 0
$ws.Range("O516").Value = 66
$ws.Range("P516").Value = 8
$ws.Range("Q516").Value = 17
$ws.Range("R516").Value = 16
$ws.Range("S516").Value = 61
$ws.Range("T516").Value = 74

$ws.Range("B517").Value = 516
$ws.Range("C517").Value = 20471
$ws.Range("D517").Value = 581
$ws.Range("E517").Value = 85
$ws.Range("F517").Value = 0.02838161301
$ws.Range("G517").Value = 19746
$ws.Range("H517").Value = 27040
$ws.Range("I517").Value = 47511
$ws.Range("J517").Value = 0
$ws.Range("K517").Value = 0
$ws.Range("L517").Value = 0
$ws.Range("M517").Value = 66
$ws.Range("N517").Value = 0
$ws.Range("O517").Value = 66
$ws.Range("P517").Value = 11
$ws.Range("Q517").Value = 16
$ws.Range("R517").Value = 16
$ws.Range("S517").Value = 61
$ws.Range("T517").Value = 74

# --- Step 2: write column A (the date strings) as TEXT, not auto-converted dates ---
# Temporarily format as Text so Excel does not coerce "yyyy/mm/dd"-looking
# strings into date serials, then clear the temporary formatting again so no
# extra style survives in the saved file.
$dateRng = $ws.Range("A511:A517")
$dateRng.NumberFormat = "@"
$ws.Range("A511").Value = "2021/08/23"
$ws.Range("A512").Value = "2021/08/24"
$ws.Range("A513").Value = "2021/08/25"
$ws.Range("A514").Value = "2021/08/26"
$ws.Range("A515").Value = "2021/08/27"
$ws.Range("A516").Value = "2021/08/28"
$ws.Range("A517").Value = "2021/08/29"
$dateRng.ClearFormats()

# --- Step 3: fix the floating point noise on previously-stored ratios in column F ---
$ws.Range("F49").Value = 0.03846153846
$ws.Range("F50").Value = 0.03773584906
$ws.Range("F51").Value = 0.05263157895
$ws.Range("F52").Value = 0.05263157895
$ws.Range("F54").Value = 0.06060606061
$ws.Range("F89").Value = 0.03345724907
$ws.Range("F140").Value = 0.02951838426
$ws.Range("F406").Value = 0.02810322989
$ws.Range("F449").Value = 0.956182004
$ws.Range("F450").Value = 0.02871098404
$ws.Range("F455").Value = 0.02837853301
$ws.Range("F506").Value = 0.02840126009
